# Devlog edit: "added ServerCommands for info about hitting other players"
#
# 1) Extend the "Added PlayerShoot Script" bullet with a trailing sentence
#    describing what the script does.
# 2) Split the "added updating of gunBarrel rotation for multiplayer" bullet:
#      - capitalize the leading "a" (typed as two runs: "A" + rest)
#      - move the trailing bookmark onto its own new "Added " bullet
#      - move the page-break onto its own (non-list) paragraph
#      - add two further new empty paragraphs (plain, then list-style)
#
# Both edits are implemented with Range.InsertXML so that the untouched
# runs/paragraphs around the edit points keep their original content
# exactly (InsertXML replaces only the exact range it is called on).

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1) "Added PlayerShoot Script" -> add a new trailing run in the same
#    paragraph explaining what the script does.
# ---------------------------------------------------------------------
$f1 = $d.Content
$null = $f1.Find.Execute("PlayerShoot Script")
$r1 = $d.Range($f1.Start, $f1.End)

$body1 = '<w:p>' + `
  '<w:r w:rsidR="005434A9"><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>PlayerShoot Script</w:t></w:r>' + `
  '<w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve"> which does a RayCast in the direction the player is aiming and then tells if we hit a RemotePlayer(every player which is not our clients)</w:t></w:r>' + `
  '</w:p>'

$null = $r1.InsertXML($pkgHeader + $body1 + $pkgFooter)

# ---------------------------------------------------------------------
# 2) Split the "added updating of gunBarrel..." paragraph into five
#    paragraphs.
# ---------------------------------------------------------------------
$f2 = $d.Content
$null = $f2.Find.Execute("added updating of gunBarrel rotation for multiplayer")
$p2 = $f2.Paragraphs(1).Range
$r2 = $d.Range($p2.Start, $p2.End)

$body2 = `
  '<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>A</w:t></w:r>' + `
    '<w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t>dded updating of gunBarrel rotation for multiplayer</w:t></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:t xml:space="preserve">Added </w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:lang w:val="de-CH"/></w:rPr><w:br w:type="page"/></w:r>' + `
  '</w:p>' + `
  '<w:p><w:pPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p>' + `
  '<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr><w:lang w:val="de-CH"/></w:rPr></w:pPr></w:p>'

$null = $r2.InsertXML($pkgHeader + $body2 + $pkgFooter)

Write-Output "done"
